# Updated symbol list on Thu Feb  9 16:44:03 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) cells on Sheet1 are
# stored as plain text (e.g. "318.58", "-2.78%") rather than numbers, so
# each value below is written back as text too (NumberFormat "@" keeps
# Excel from auto-converting the numeric/percent-looking strings into
# real numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{ Cell = "D2"; Value = "318.16" },
    @{ Cell = "E2"; Value = "-3.26%" },
    @{ Cell = "D3"; Value = "42.04" },
    @{ Cell = "E3"; Value = "-4.55%" },
    @{ Cell = "D4"; Value = "5.176" },
    @{ Cell = "E4"; Value = "-4.30%" },
    @{ Cell = "D5"; Value = "0.08125" },
    @{ Cell = "E5"; Value = "-2.26%" },
    @{ Cell = "E6"; Value = "-1.82%" },
    @{ Cell = "D7"; Value = "1.759" },
    @{ Cell = "E7"; Value = "-9.04%" },
    @{ Cell = "D8"; Value = "0.9331" },
    @{ Cell = "E8"; Value = "-3.77%" },
    @{ Cell = "D9"; Value = "0.1131" },
    @{ Cell = "E9"; Value = "1.92%" },
    @{ Cell = "D10"; Value = "0.1857" },
    @{ Cell = "E10"; Value = "-1.59%" },
    @{ Cell = "D11"; Value = "0.04688" },
    @{ Cell = "E11"; Value = "-0.09%" },
    @{ Cell = "D12"; Value = "0.09271" },
    @{ Cell = "E12"; Value = "-4.82%" },
    @{ Cell = "D13"; Value = "7.415" },
    @{ Cell = "E13"; Value = "-19.22%" },
    @{ Cell = "D14"; Value = "0.1054" },
    @{ Cell = "E14"; Value = "-0.48%" },
    @{ Cell = "D15"; Value = "0.001277" },
    @{ Cell = "E15"; Value = "-2.47%" },
    @{ Cell = "D16"; Value = "0.005756" },
    @{ Cell = "E16"; Value = "-4.35%" },
    @{ Cell = "D17"; Value = "3.354" },
    @{ Cell = "E17"; Value = "-0.95%" },
    @{ Cell = "D18"; Value = "2.551" },
    @{ Cell = "E18"; Value = "0.94%" },
    @{ Cell = "D19"; Value = "0.3388" },
    @{ Cell = "E19"; Value = "1.92%" },
    @{ Cell = "D20"; Value = "0.1396" },
    @{ Cell = "E20"; Value = "1.03%" },
    @{ Cell = "D21"; Value = "0.2602" },
    @{ Cell = "E21"; Value = "-4.35%" },
    @{ Cell = "D22"; Value = "0.04182" },
    @{ Cell = "E22"; Value = "-0.24%" },
    @{ Cell = "D23"; Value = "0.001246" },
    @{ Cell = "E23"; Value = "-4.40%" },
    @{ Cell = "D24"; Value = "0.004276" },
    @{ Cell = "E24"; Value = "-3.34%" },
    @{ Cell = "D25"; Value = "0.0001225" },
    @{ Cell = "E25"; Value = "-5.97%" },
    @{ Cell = "D26"; Value = "0.0002995" },
    @{ Cell = "E26"; Value = "0.27%" },
    @{ Cell = "D38"; Value = "0.02583" },
    @{ Cell = "E38"; Value = "-2.46%" },
    @{ Cell = "D39"; Value = "0.05492" },
    @{ Cell = "E39"; Value = "-1.85%" },
    @{ Cell = "D40"; Value = "0.008138" },
    @{ Cell = "E40"; Value = "3.47%" },
    @{ Cell = "D41"; Value = "0.1394" },
    @{ Cell = "E41"; Value = "-0.83%" },
    @{ Cell = "D42"; Value = "0.006542" },
    @{ Cell = "E42"; Value = "-11.54%" },
    @{ Cell = "D43"; Value = "0.002088" },
    @{ Cell = "E43"; Value = "-1.48%" },
    @{ Cell = "D44"; Value = "0.007624" },
    @{ Cell = "E44"; Value = "-11.72%" },
    @{ Cell = "D45"; Value = "0.3475" },
    @{ Cell = "E45"; Value = "3.17%" },
    @{ Cell = "D46"; Value = "0.00006811" },
    @{ Cell = "E46"; Value = "-0.29%" },
    @{ Cell = "E47"; Value = "0.17%" },
    @{ Cell = "D48"; Value = "0.003383" },
    @{ Cell = "E48"; Value = "-3.76%" },
    @{ Cell = "D49"; Value = "0.004122" },
    @{ Cell = "E49"; Value = "16.55%" },
    @{ Cell = "E50"; Value = "0.17%" },
    @{ Cell = "D51"; Value = "0.0002008" },
    @{ Cell = "E51"; Value = "0.17%" }
)

foreach ($update in $cellUpdates) {
    $rng = $ws.Range($update.Cell)
    # Force text storage so Excel doesn't reinterpret "318.16" / "-3.26%"
    # style strings as a number / percentage.
    $rng.NumberFormat = "@"
    $rng.Value = $update.Value
}
